# Reorder "Recorded By" (column G) entries so that "System" is moved
# from the front of the comma-separated list to the end, e.g.
#   "System, backup@backdoor.com" -> "backup@backdoor.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ($text -ne $null -and $text -ne "") {
        if ($text.StartsWith("System, ")) {
            $parts = $text.Split(", ")
            $rest = $parts[1..($parts.Length - 1)]
            $newText = ($rest -join ", ") + ", System"
            $cell.Value = $newText
        }
    }
}
